$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 45, pushing the existing
# rows 45-46 down to become rows 47-48 (their content is preserved as-is
# by the insert/shift, matching the target diff).
$ws.Range("A45:A46").EntireRow.Insert()

# New row 45: updated "Primera" quality entry for the newer date (44448)
$ws.Range("A45").Value = 9
$ws.Range("B45").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 44448
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100102
$ws.Range("H45").Value = "Cítricos"
$ws.Range("I45").Value = 100102006
$ws.Range("J45").Value = "Pomelo"
$ws.Range("K45").Value = "Start Ruby"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 140
$ws.Range("N45").Value = 7000
$ws.Range("O45").Value = 7000
$ws.Range("P45").Value = 7000
$ws.Range("Q45").Value = "$/caja 14 kilos granel"
$ws.Range("R45").Value = "Región Metropolitana"
$ws.Range("S45").Value = 500
$ws.Range("T45").Value = 14

# New row 46: brand-new "Segunda" quality entry for the newer date (44448)
$ws.Range("A46").Value = 9
$ws.Range("B46").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C46").Value = "Metropolitana"
$ws.Range("D46").Value = 44448
$ws.Range("E46").Value = 13
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100102
$ws.Range("H46").Value = "Cítricos"
$ws.Range("I46").Value = 100102006
$ws.Range("J46").Value = "Pomelo"
$ws.Range("K46").Value = "Start Ruby"
$ws.Range("L46").Value = "Segunda"
$ws.Range("M46").Value = 220
$ws.Range("N46").Value = 5000
$ws.Range("O46").Value = 5000
$ws.Range("P46").Value = 5000
$ws.Range("Q46").Value = "$/caja 14 kilos granel"
$ws.Range("R46").Value = "Región Metropolitana"
$ws.Range("S46").Value = 357
$ws.Range("T46").Value = 14
